$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Administrator, Miss Dina Nasr") {
        $cell.Value2 = "Miss Dina Nasr, Administrator"
    }
}
